$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 12-13 (existing rows 12-21 shift down to 14-23).
# Excel's insert behaviour copies per-cell formatting down from the row
# above (row 11), which already carries the s11/s5/s3/s9/s6 style pattern
# used by this "two sub-technique" block layout.
$ws.Rows("12:13").Insert()
$ws.Rows("12:13").RowHeight = 19

# --- Row 12: Reverse Pointer technique ---
$ws.Range("A12").Value = "Arrays"
$ws.Range("B12").Value = "Reverse Pointer"
$ws.Range("C12").Value = "https://leetcode.com/problems/reverse-only-letters"
$ws.Range("D12").Value = "ReverseOnlyLetters_ReversePointer"
$ws.Range("E12").Value = "Reverse Pointer"

# --- Row 13: Stack technique ---
$ws.Range("A13").Value = "Arrays"
$ws.Range("B13").Value = "Reverse Pointer"
$ws.Range("C13").Value = "https://leetcode.com/problems/reverse-only-letters"
$ws.Range("D13").Value = "ReverseOnlyLetters_Stack"
$ws.Range("E13").Value = "Stack (Pop, Push)"

# Notes for row 13 written first (matches original authoring/shared-string order)
$ws.Range("F13").Value = "Use Stack  Technique (Collect the letters of 'S' separately into a stack, so that popping the stack reverses the letters)"
$ws.Range("F13").Font.Size = 14
$ws.Range("F13").Font.Name = "Arial (Body)"
$ws.Range("F13").Characters(4, 7).Font.Color = 255
$ws.Range("F13").Characters(11, 111).Font.Color = 0

# Notes for row 12 written last
$ws.Range("F12").Value = "Use Reverse Pointer Technique: You need two loop, the first loop will check if the left character can be reversed, if so then we need another loop that will check if right character can be reversed   (int j = S.length() - 1;   while (!Character.isLetter(S.charAt(j))) j--)  "
$ws.Range("F12").Font.Size = 14
$ws.Range("F12").Font.Name = "Arial (Body)"
$ws.Range("F12").Characters(4, 16).Font.Color = 255
$ws.Range("F12").Characters(20, 180).Font.Color = 0
$ws.Range("F12").Characters(200, 75).Font.Color = 255

# Hyperlinks for the two new rows (re-apply the plain "Hyperlink" cell
# style afterwards so the cell keeps reusing the workbook's existing
# Hyperlink style index instead of leaving a freshly-applied variant).
$ws.Hyperlinks.Add($ws.Range("C13"), "https://leetcode.com/problems/reverse-only-letters") | Out-Null
$ws.Range("C13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C12"), "https://leetcode.com/problems/reverse-only-letters") | Out-Null
$ws.Range("C12").Style = "Hyperlink"

# Restore the view state: scrolled so column F is visible starting at row 3,
# with F12 as the active selection.
$ws.Range("F12").Select()
